$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compOptGebradd")
$co = $ws.ChartObjects("Chart 5")
$ch = $co.Chart
Write-Host "ChartType: $($ch.ChartType)"
Write-Host "HasTitle: $($ch.HasTitle)"
$ax1 = $ch.Axes(1)
Write-Host "Axis1 HasTitle: $($ax1.HasTitle)"
if ($ax1.HasTitle) {
    Write-Host "Axis1 Title: $($ax1.AxisTitle.Text)"
}
